$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-01 Sunday" "2025-06-02 Monday"

Replace-Text "66÷9=" "18÷4="
Replace-Text "55÷5=" "62÷3="
Replace-Text "22÷6=" "85÷8="
Replace-Text "53÷6=" "80÷6="
Replace-Text "41÷8=" "37÷3="

Replace-Text "82÷3=" "98÷6="
Replace-Text "39÷6=" "85÷3="
Replace-Text "18÷5=" "75÷4="
Replace-Text "80÷8=" "82÷3="
Replace-Text "67÷5=" "69÷3="

Replace-Text "78÷9=" "15÷2="
Replace-Text "10÷2=" "52÷4="
Replace-Text "36÷6=" "21÷4="
Replace-Text "23÷4=" "74÷2="
Replace-Text "63÷3=" "81÷8="

Replace-Text "68÷6=" "39÷7="
Replace-Text "63÷4=" "13÷6="
Replace-Text "62÷8=" "23÷9="
Replace-Text "11÷9=" "56÷4="
Replace-Text "51÷4=" "94÷9="

Replace-Text "82÷2=" "17÷6="
Replace-Text "90÷7=" "69÷6="
Replace-Text "12÷6=" "93÷5="
Replace-Text "90÷9=" "35÷5="
Replace-Text "54÷3=" "28÷2="
